$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column headers: PREREQ_COURSES -> Prerequisites, COREQ_COURSES -> CoRequisites
$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Update the selection to A1:I1 (header row)
$ws.Range("A1:I1").Select()
